# Update the "Rules" worksheet: change the greeting text in E8 from
# "Good Morning" to "GIT UPDATE" and leave the last-selected cell as E8
# (matching the state Excel would persist after a user typed into E8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
